# Backlog.xlsx update — adds new backlog items (länkar/produktkategorier)
# and swaps the "Sprint"/"Kommentar" columns (E <-> F), widening the
# comment column to fit the new longer notes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Swap the E/F header labels ("Kommentar" <-> "Sprint") ---------
$ws.Range("E1").Value = "Sprint"
$ws.Range("F1").Value = "Kommentar"

# --- 2. Move the existing comment/sprint notes to their new columns ---
# Row 10 had its note in E, row 11 had its note in F; after the column
# swap they live in F and E respectively.
$ws.Range("F10").Value = $ws.Range("E10").Value2
$ws.Range("E10").ClearContents()

$ws.Range("E11").Value = $ws.Range("F11").Value2
$ws.Range("F11").ClearContents()

# --- 3. Insert three new backlog rows before the old row 13 -----------
$ws.Rows(13).Resize(3).Insert()

# New row 13: Produktkategorier (Databas / Hög / Dålig) + comment
$ws.Range("A13").Value = "Produktkategorier"
$ws.Range("B13").Value = "Databas"
$ws.Range("C13").Value = "Hög"
$ws.Range("D13").Style = "Dålig"
$ws.Range("F13").Value = "Istället för bara strängar i Products så har vi foreign key till en tabell."

# New row 14: Gemensam meny för hemsidan (Front end / Hög / Dålig)
$ws.Range("A14").Value = "Gemensam meny för hemsidan"
$ws.Range("B14").Value = "Front end"
$ws.Range("C14").Value = "Hög"
$ws.Range("D14").Style = "Dålig"

# New row 15: Produktlänkar från databasen (Front end / Medel / Neutral) + comment
$ws.Range("A15").Value = "Produktlänkar från databasen"
$ws.Range("B15").Value = "Front end"
$ws.Range("C15").Value = "Medel"
$ws.Range("D15").Style = "Neutral"
$ws.Range("F15").Value = "Hämta produkter från databasen och genera klickbara länkar till produkterna."

# --- 4. Column widths: split D/E apart, narrow E, widen F -------------
$ws.Columns("E").ColumnWidth = 8.166666666666666
$ws.Columns("F").ColumnWidth = 76.66666666666667

# --- 5. Restore the active selection to the newly added comment cell --
$ws.Range("F15").Select()
